# "Generate Report for Handback"
# The e4e40e53-... entry has now been handed back (in sync with en-US),
# so update its status/date/error fields on every sheet, then re-sort each
# table by file name ascending (which moves the e4e40e53 row to the top),
# and finally refresh the hyperlink display text so it matches the file
# now shown in each row.

$wb = $excel.ActiveWorkbook

$ovw   = $wb.Worksheets.Item("Overview")
$zhcn  = $wb.Worksheets.Item("zh-cn")
$dede  = $wb.Worksheets.Item("de-de")

$handedBack = "Handed back: in sync with en-US"

# --- Overview sheet: row 4 is the e4e40e53 entry before sorting ---
$ovw.Range("E4").Value = $handedBack
$ovw.Range("F4").Value = $handedBack

# --- zh-cn sheet: row 4 is the e4e40e53 entry before sorting ---
$zhcn.Range("C4").Value = $handedBack
$zhcn.Range("K4").Value = "2016-09-08 05:25:48"
$zhcn.Range("P4").Value = ""

# --- de-de sheet: row 4 is the e4e40e53 entry before sorting ---
$dede.Range("C4").Value = $handedBack
$dede.Range("K4").Value = "2016-09-08 05:25:57"
$dede.Range("P4").Value = ""

# --- Re-sort every table ascending by its first (file name) column ---
$ovwTable = $ovw.ListObjects.Item(1)
$ovwTable.Sort.SortFields.Clear()
$ovwTable.Sort.SortFields.Add($ovw.Range("A2:A4"))
$ovwTable.Sort.Header = 1
$ovwTable.Sort.Apply()

$zhcnTable = $zhcn.ListObjects.Item(1)
$zhcnTable.Sort.SortFields.Clear()
$zhcnTable.Sort.SortFields.Add($zhcn.Range("A2:A4"))
$zhcnTable.Sort.Header = 1
$zhcnTable.Sort.Apply()

$dedeTable = $dede.ListObjects.Item(1)
$dedeTable.Sort.SortFields.Clear()
$dedeTable.Sort.SortFields.Add($dede.Range("A2:A4"))
$dedeTable.Sort.Header = 1
$dedeTable.Sort.Apply()

# --- Refresh hyperlink display text so it matches the row's new content ---
foreach ($h in $ovw.Hyperlinks) {
    $h.TextToDisplay = $h.Range.Text
}

foreach ($h in $zhcn.Hyperlinks) {
    $h.TextToDisplay = $h.Range.Text
}

foreach ($h in $dede.Hyperlinks) {
    $h.TextToDisplay = $h.Range.Text
}
